$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 0.02502539014123054
$ws.Range("F2").Value = "norm_coldread_gaze_wpm_median"

$ws.Range("E3").Value = 0.003128925285112487
$ws.Range("F3").Value = "norm_coldread_saccade_regression_rate_%"

$ws.Range("E4").Value = 0.03959755131507949
$ws.Range("F4").Value = "norm_qa_coverage_line_%"
